$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target row is overwritten with the full (Fecha/Calidad/Volumen/
# Precio.../Unidad/Origen/Precio-Kg/Kg-unidad) record that used to live on
# a different row of the same sheet (weekly re-shuffle of the price rows;
# row 7 is left untouched). Map: destination row -> source row (values read
# from the ORIGINAL, pre-edit sheet state).
$rowMap = @{
    2  = 5
    3  = 6
    4  = 11
    5  = 12
    6  = 8
    8  = 3
    9  = 4
    10 = 9
    11 = 10
    12 = 13
    13 = 2
}

# Columns copied as part of each record: D, L, M, N, O, P, Q, R, S, T.
$cols = @(4, 12, 13, 14, 15, 16, 17, 18, 19, 20)

# Snapshot the original values for every row/column involved before any
# writes happen, so later writes never read already-mutated cells.
$snapshot = @{}
foreach ($r in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($r)) {
        $rec = @{}
        foreach ($c in $cols) {
            $rec[$c] = $ws.Cells.Item($r, $c).Value2
        }
        $snapshot[$r] = $rec
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rec = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($destRow, $c).Value2 = $rec[$c]
    }
}
